# Apply crypto price/volume updates from the Coinranking scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force the cell to keep a literal text value (matches the inline-string
    # cells already used in this sheet) instead of Excel auto-converting
    # numeric-looking strings (e.g. "246.17", "1.00") into real numbers.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '37.013.63'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '2.055.47'
$ws.Range('E3').Value = '  +0.32%  '
$ws.Range('E4').Value = '  -0.06%  '
Set-TextValue $ws.Range('D5') '246.17'
$ws.Range('E5').Value = '  -1.13%  '
$ws.Range('E6').Value = '  -1.55%  '
Set-TextValue $ws.Range('D7') '58.55'
$ws.Range('E7').Value = '  -1.02%  '
$ws.Range('E8').Value = '  -0.03%  '
Set-TextValue $ws.Range('D9') '0.377'
$ws.Range('E9').Value = '  -2.60%  '
Set-TextValue $ws.Range('D10') '0.0775'
$ws.Range('E10').Value = '  -1.95%  '
$ws.Range('E11').Value = '  +2.35%  '
Set-TextValue $ws.Range('D12') '15.45'
$ws.Range('E12').Value = '  -3.10%  '
Set-TextValue $ws.Range('D13') '0.884'
$ws.Range('E13').Value = '  +6.31%  '
$ws.Range('D14').Value = '2.355.88'
$ws.Range('E14').Value = '  +0.22%  '
Set-TextValue $ws.Range('D15') '5.70'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = '2.090.33'
$ws.Range('E16').Value = '  +1.95%  '
Set-TextValue $ws.Range('D17') '18.16'
$ws.Range('E17').Value = '  -2.60%  '
$ws.Range('D18').Value = '36.982.74'
$ws.Range('E18').Value = '  -0.19%  '
Set-TextValue $ws.Range('D19') '73.80'
$ws.Range('E19').Value = '  -2.11%  '
$ws.Range('D20').Value = '0.0₃0890'
$ws.Range('E20').Value = '  -1.31%  '
Set-TextValue $ws.Range('D21') '5.44'
$ws.Range('E21').Value = '  +0.39%  '
Set-TextValue $ws.Range('D22') '238.08'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  +1.60%  '
Set-TextValue $ws.Range('D25') '10.18'
$ws.Range('E25').Value = '  +6.58%  '
Set-TextValue $ws.Range('D26') '169.79'
$ws.Range('E26').Value = '  +0.69%  '
$ws.Range('E27').Value = '  -2.25%  '
Set-TextValue $ws.Range('D28') '20.11'
$ws.Range('E28').Value = '  +0.14%  '
Set-TextValue $ws.Range('D29') '5.50'
$ws.Range('E29').Value = '  +14.96%  '
$ws.Range('E30').Value = '  -1.27%  '
$ws.Range('E31').Value = '  -1.49%  '
$ws.Range('E32').Value = '  +3.36%  '
$ws.Range('E33').Value = '  -1.30%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D34') '2.34'
$ws.Range('E34').Value = '  +5.42%  '
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D35') '1.00'
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('E36').Value = '  +5.44%  '
Set-TextValue $ws.Range('D37') '0.0846'
$ws.Range('E37').Value = '  -5.16%  '
$ws.Range('E38').Value = '  -0.16%  '
Set-TextValue $ws.Range('D39') '5.23'
$ws.Range('E39').Value = '  +2.92%  '
$ws.Range('E40').Value = '  -1.70%  '
$ws.Range('E41').Value = '  +0.26%  '
$ws.Range('E42').Value = '  +2.10%  '
Set-TextValue $ws.Range('D43') '0.0959'
$ws.Range('E43').Value = '  -10.52%  '
Set-TextValue $ws.Range('D44') '97.35'
$ws.Range('E44').Value = '  +0.57%  '
Set-TextValue $ws.Range('D45') '16.96'
$ws.Range('E45').Value = '  -4.04%  '
$ws.Range('D46').Value = '1.301.20'
$ws.Range('E46').Value = '  +0.95%  '
Set-TextValue $ws.Range('D47') '2.36'
$ws.Range('E47').Value = '  -5.52%  '
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').Value = '2.245.31'
$ws.Range('E50').Value = '  +0.47%  '
$ws.Range('E51').Value = '  +2.57%  '
